$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "ink, indelible visible ink"
$ws.Range("B5").Value = "small , mountainous state"
$ws.Range("B8").Value = "country 's elections, elections, recent elections, many elections"
$ws.Range("B9").Value = "greatest part, part"
$ws.Range("B10").Value = "petition drive, drive"
$ws.Range("B11").Value = "actual technology, this new technology"
$ws.Range("B18").Value = "use, improper use"
$ws.Range("B44").Value = "other common type, this type"
